$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.192.89'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.586.29'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.809.55'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '1.598.34'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '26.199.75'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '212.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0496'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.86%  '
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').Value = '1.333.63'
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.581'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.79%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.819'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  +3.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.951'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -13.61%  '
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = '1.721.34'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('E50').Value = '  -2.07%  '

Write-Host "Applied 67 cell updates"
